# Commit: "Updated RMD and PPX"
#
# The deck had a stray duplicate of the "Top Three Causes of Attrition"
# slide (SlideID 272) sitting right before the closing "Root Mean square
# error ..." slide. Remove that duplicate slide, leaving the original
# "Top Three Causes of Attrition" slide (earlier in the deck, right
# before "Prediction using KNN") untouched.
#
# Locate the slide by its stable SlideID (272, matching sldId="272" /
# cId="1911460457" in the deck's change history) rather than a
# hard-coded index, so the edit is unambiguous even if slides are ever
# reordered.

$p = $ppt.ActivePresentation

$targetSlideId = 272

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    if ($s.SlideID -eq $targetSlideId) {
        $s.Delete()
        break
    }
}
